$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("g3.6")

# Update existing values in row 15 (2023)
$ws.Range("B15").Value = 14.29120608080826
$ws.Range("C15").Value = 8.064344405632822

# Add new row 16 (2024)
$ws.Range("A16").Value = 2024
$ws.Range("B16").Value = 16.66347210408774
$ws.Range("C16").Value = 9.443802699703051
$ws.Range("D16").Value = 22.62152262886292
